$d = $word.ActiveDocument

# --- 1. Insert two new paragraphs after the "pcb spark gap" paragraph ---
# Locate the paragraph that ends with the spark-gap / R3 sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*A pcb spark gap is made and labelled pad PD1*") {
        $target = $p
        break
    }
}

# Insert a blank paragraph right after it (mirrors the existing blank-line
# spacing pattern used throughout the document).
$target.Range.InsertParagraphAfter()

# The blank paragraph we just created is now the paragraph immediately
# following $target; insert a second new paragraph after that one to hold
# the new sentence.
$blankPara = $target.Next()
$blankPara.Range.InsertParagraphAfter()

# That second new paragraph is the one right after $blankPara; fill it in
# with the new text.
$textPara = $blankPara.Next()
$textPara.Range.InsertAfter("Thru-hole pads located near the pcb spark gap for mounting a gas discharge type spark gap.")

# --- 2. Update the cached "PAGE" field result in the header from 8 to 9 ---
$hdr = $d.Sections(1).Headers(1)
$null = $hdr.Range.Find.Execute("8", $true, $false, $false, $false, $false, $true, 1, $false, "9", 2)
